$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.567.31'
$ws.Range('E2').Value = '  -1.05%  '
$ws.Range('D3').Value = '2.562.81'
$ws.Range('E3').Value = '  -0.23%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '578.07'
$ws.Range('E5').Value = '  -1.18%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.10'
$ws.Range('E6').Value = '  -4.07%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -0.20%  '
$ws.Range('E9').Value = '  -3.37%  '
$ws.Range('E10').Value = '  -1.25%  '
$ws.Range('E11').Value = '  -0.75%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.349'
$ws.Range('E12').Value = '  -2.30%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '26.72'
$ws.Range('E13').Value = '  -4.57%  '
$ws.Range('D14').Value = '3.022.46'
$ws.Range('E14').Value = '  -0.21%  '
$ws.Range('D15').Value = '62.472.40'
$ws.Range('E15').Value = '  -1.03%  '
$ws.Range('E16').Value = '  -3.19%  '
$ws.Range('D17').Value = '2.554.65'
$ws.Range('E17').Value = '  -0.32%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.07'
$ws.Range('E18').Value = '  -3.46%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '337.83'
$ws.Range('E19').Value = '  -1.13%  '
$ws.Range('E20').Value = '  -2.00%  '
$ws.Range('E21').Value = '  -3.50%  '
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.10'
$ws.Range('E23').Value = '  +1.35%  '
$ws.Range('E24').Value = '  -5.24%  '
$ws.Range('E25').Value = '  -4.28%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.50'
$ws.Range('E26').Value = '  +0.63%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('E28').Value = '  -5.09%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.15'
$ws.Range('E29').Value = '  -5.07%  '
$ws.Range('E30').Value = '  -2.92%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '453.76'
$ws.Range('E31').Value = '  +2.68%  '
$ws.Range('D32').Value = '0.0₃0794'
$ws.Range('E32').Value = '  -4.72%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '176.74'
$ws.Range('E33').Value = '  -0.50%  '
$ws.Range('E34').Value = '  +0.22%  '
$ws.Range('E35').Value = '  +0.09%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.398'
$ws.Range('E36').Value = '  -2.72%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.79'
$ws.Range('E37').Value = '  -3.00%  '
$ws.Range('E38').Value = '  -2.77%  '
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('E40').Value = '  -5.10%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '40.40'
$ws.Range('E41').Value = '  +1.53%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '157.41'
$ws.Range('E42').Value = '  +3.22%  '
$ws.Range('E43').Value = '  -4.30%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.628'
$ws.Range('E44').Value = '  +2.85%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '20.78'
$ws.Range('E45').Value = '  -3.68%  '
$ws.Range('E46').Value = '  -4.80%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0955'
$ws.Range('E47').Value = '  -2.53%  '
$ws.Range('E48').Value = '  -4.24%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '17.89'
$ws.Range('E49').Value = '  -3.51%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '11.41'
$ws.Range('E50').Value = '  +0.35%  '
$ws.Range('E51').Value = '  -5.95%  '
